# Entsoe Actual Production Solar update:
# - Shift all timestamp values (column A, rows 2-97) forward by 10 days.
# - Update the actual production values (column B, rows 23-39) with the
#   newly fetched data for the shifted date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2 through 97) forward by 10 days.
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2 + 10
}

# New "Actual Production (MW)" values for the shifted day (column B).
$newValues = @{
    23 = 6
    24 = 19
    25 = 36
    26 = 69
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
